# Fix two typos on the "Költésgvetés" (budget) sheet and update the
# selected cell, matching the upstream re-upload of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Költésgvetés")

# A5: "Interaktí tábla" -> "Interaktív tábla"
$ws.Range("A5").Value = "Interaktív tábla"

# A12: "Accespoint" -> "Accesspoint"
$ws.Range("A12").Value = "Accesspoint"

# Move the sheet's active selection from E16 to A13.
$ws.Range("A13").Select()
